# "Added last minute updates"
# - Rename the placeholder ID token in the first paragraph from
#   **ID__AFFARS_5330_topic_2__ID** to **ID__AFFARS_5330_201_5__ID**,
#   merging away the trailing space-only run that followed it.
# - Give that same paragraph the paragraph border (5-twip space on all
#   sides) and the 225-twip left indent that the other body paragraphs
#   already use.

$d = $word.ActiveDocument

# Replace the old placeholder text (plus the trailing space that lived in
# its own run) with the new placeholder text - this both updates the text
# and drops the now-unwanted trailing-space run.
$d.Content.Find.Execute(
    "AFFARS_5330_topic_2__ID** ", $true, $false, $false, $false, $false,
    $true, 1, $false, "AFFARS_5330_201_5__ID**", 2)

# First paragraph of the document: add the paragraph border and update
# the left indent to match the rest of the body paragraphs.
$p = $d.Paragraphs(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25
